# Regenerate merged AHB files
#
# 1) Rename the "old" / "new" header-suffix columns to the new version tags
#    (FV2410 / FV2504) used by the regenerated AHB diff export.
# 2) Turn the data range into a proper Excel Table ("Table1").
# 3) Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header labels -------------------------------------------------
# Columns A:J were suffixed "_old" -> now "_FV2410"
# Columns L:U were suffixed "_new" -> now "_FV2504"
# Column K ("diff") is left untouched.
$headerRange = $ws.Range("A1:U1")
$null = $headerRange.Replace("_old", "_FV2410", 2)
$null = $headerRange.Replace("_new", "_FV2504", 2)

# --- 2) Convert the data range A1:U62 into an Excel Table --------------------
$dataRange = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row -------------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
